$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the stray <w:bookmarkStart/><w:bookmarkEnd/> ("_GoBack") that sat
#    just before the "Schedule by Week:" paragraph.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2) After the "No sound designers" bullet, add four new sub-bullets:
#      Engine            (ilvl 0)
#      Sound              (ilvl 1)
#      Console input      (ilvl 1)
#      Tweaking           (ilvl 1)  <- also carries the relocated _GoBack bookmark
# ---------------------------------------------------------------------------
$targetNum = 0
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -like "*No sound designers*") {
        $targetNum = $i
    }
}

$target = $d.Paragraphs.Item($targetNum)
$target.Range.InsertParagraphAfter()
$target.Range.InsertParagraphAfter()
$target.Range.InsertParagraphAfter()
$target.Range.InsertParagraphAfter()

$texts  = @("Engine", "Sound", "Console input", "Tweaking")
$levels = @(1, 2, 2, 2)

for ($k = 1; $k -le 4; $k++) {
    $p = $d.Paragraphs.Item($targetNum + $k)
    $p.Range.Text = $texts[$k - 1]
    $p.Range.ListFormat.ListLevelNumber = $levels[$k - 1]
}

# ---------------------------------------------------------------------------
# Re-create the "_GoBack" bookmark right after the "Tweaking" run's text
# (matching the original placement semantics: bookmarkStart/bookmarkEnd both
# appear immediately after the run, not wrapping it). A directly-collapsed
# Range positioned at "end-of-paragraph-text" is unreliable, so a temporary
# marker character is appended, the (now safely mid-paragraph) boundary is
# bookmarked, and the marker is deleted again -- the bookmark collapses in
# place and survives the deletion.
# ---------------------------------------------------------------------------
$lastp = $d.Paragraphs.Item($targetNum + 4)
$boundary = $lastp.Range.End - 1
$d.Range($boundary, $boundary).InsertAfter("X")
$d.Bookmarks.Add("_GoBack", $d.Range($boundary, $boundary))
$d.Range($boundary, $boundary + 1).Delete()
